$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns keep their text formatting
# so numeric-looking strings (e.g. "292.73", "0.999") are not converted
# into real numbers by Excel, matching the original inline-string cells.
$ws.Range("D2:E51").NumberFormat = "@"

$updates = @{
    "D2" = "39.848.66"
    "E2" = "  -0.27%  "
    "D3" = "2.202.57"
    "E3" = "  -1.10%  "
    "D4" = "0.999"
    "E4" = "  -0.06%  "
    "D5" = "292.73"
    "E5" = "  -0.59%  "
    "D6" = "86.86"
    "E7" = "  -1.16%  "
    "E8" = "  -0.03%  "
    "E9" = "  -0.09%  "
    "B10" = "Avalanche"
    "C10" = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
    "D10" = "29.80"
    "E10" = "  -3.95%  "
    "B11" = "Dogecoin"
    "C11" = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
    "D11" = "0.0774"
    "E11" = "  -1.88%  "
    "D12" = "49.43"
    "E12" = "  +5.40%  "
    "E13" = "  +2.57%  "
    "D14" = "6.41"
    "E14" = "  -0.35%  "
    "D15" = "2.544.76"
    "E15" = "  -0.99%  "
    "D16" = "13.67"
    "E16" = "  -2.65%  "
    "D17" = "2.216.28"
    "E17" = "  -5.35%  "
    "E18" = "  -0.54%  "
    "D19" = "39.737.97"
    "E19" = "  -0.32%  "
    "E20" = "  -0.70%  "
    "D21" = "11.30"
    "E21" = "  +4.95%  "
    "E22" = "  -0.92%  "
    "D23" = "65.12"
    "E23" = "  -0.02%  "
    "D24" = "235.80"
    "E24" = "  +0.36%  "
    "E25" = "  +0.08%  "
    "D26" = "2.44"
    "E26" = "  -0.59%  "
    "E27" = "  -2.23%  "
    "B28" = "Toncoin"
    "C28" = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
    "D28" = "2.30"
    "E28" = "  +3.25%  "
    "B29" = "EthereumClassic"
    "C29" = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
    "D29" = "22.32"
    "E29" = "  -1.84%  "
    "D30" = "9.14"
    "E30" = "  -0.73%  "
    "D31" = "155.13"
    "E31" = "  +2.06%  "
    "D32" = "31.52"
    "E32" = "  -4.75%  "
    "D33" = "0.999"
    "E33" = "  -0.01%  "
    "D34" = "4.87"
    "E34" = "  +0.41%  "
    "D35" = "0.0707"
    "E35" = "  -1.20%  "
    "E37" = "  +4.26%  "
    "E38" = "  +0.15%  "
    "E39" = "  -2.18%  "
    "D40" = "15.38"
    "E40" = "  -5.25%  "
    "E41" = "  -1.37%  "
    "D42" = "2.120.18"
    "E42" = "  +4.54%  "
    "E43" = "  -1.99%  "
    "E44" = "  -4.44%  "
    "B45" = "VeChain"
    "C45" = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
    "D45" = "0.0266"
    "E45" = "  -1.13%  "
    "B46" = "EnergySwap"
    "C46" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "D46" = "17.66"
    "E46" = "  +9.10%  "
    "D47" = "9.58"
    "E47" = "  -3.28%  "
    "E48" = "  +2.72%  "
    "D49" = "2.418.52"
    "E49" = "  -1.07%  "
    "D50" = "1.46"
    "E50" = "  +0.40%  "
    "E51" = "  +0.72%  "
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
